# Debugged - ready to show
# Fix up the TI (Trial Inclusion/Exclusion Criteria) sheet: columns A and B on
# data rows were carrying stray header-style labels (DOMAIN, IETESTCD, ...)
# instead of the real per-row values, column C (the running criterion number)
# was blank on several rows, and column H (TIVERS) was blank/numeric-zero
# instead of the text value "1" that marks each row as belonging to version 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TI")

for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 1).Value = " "      # A: blank placeholder
    $ws.Cells.Item($r, 2).Value = "TI"     # B: DOMAIN value
    $ws.Cells.Item($r, 3).Value = [string]($r - 1)   # C: running criterion number
    $ws.Cells.Item($r, 8).Value = "1"      # H: TIVERS
}
